$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 6 (old rows 6-8 shift down to 8-10); the
# second of the two stays completely blank, leaving a gap at row 7.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# The inserted rows pick up row 5's formatting by default - strip it back off.
$ws.Rows.Item(6).ClearFormats()
$ws.Rows.Item(7).ClearFormats()
$ws.Rows.Item(7).ClearContents()

# Fill in the new row 6 data, left to right so shared-string order matches
$ws.Cells.Item(6, 1).Value = "spikes.png"
$ws.Hyperlinks.Add($ws.Cells.Item(6, 2), "https://opengameart.org/content/spikes-0")
$ws.Cells.Item(6, 3).Value = "Public Domain CC0"
$ws.Cells.Item(6, 4).Value = "No attribution"

$ws.Range("D6").Select()
